$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "cryptos" price/volume table (columns D = Price, E = Volume(1h))
# with the latest scraped values. Column D entries are prefixed with a
# leading apostrophe so Excel stores them as literal text (matching the
# source inlineStr cells) instead of auto-converting number-looking
# strings like "1.00" or "65.888.25" into numeric values.
$ws.Range("D2").Value = "'65.888.25"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "'3.327.68"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'188.85"
$ws.Range("E5").Value = "  +4.93%  "
$ws.Range("D6").Value = "'556.02"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  -0.10%  "

# Row 8 and row 9 swap places: LidoStakedEther now ranks above XRP.
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "'3.317.13"
$ws.Range("E8").Value = "  +1.46%  "
$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D9").Value = "'0.581"
$ws.Range("E9").Value = "  -1.16%  "
$ws.Range("D10").Value = "'0.180"
$ws.Range("E10").Value = "  -3.38%  "
$ws.Range("E11").Value = "  -0.62%  "
$ws.Range("D12").Value = "'46.29"
$ws.Range("E12").Value = "  -2.36%  "
$ws.Range("E13").Value = "  +1.80%  "
$ws.Range("D14").Value = "'8.56"
$ws.Range("D15").Value = "'3.856.69"
$ws.Range("E15").Value = "  +1.27%  "
$ws.Range("D16").Value = "'594.28"
$ws.Range("E16").Value = "  -6.43%  "
$ws.Range("D17").Value = "'65.960.41"
$ws.Range("E17").Value = "  +0.52%  "

# Row 18 and row 20 swap places: WrappedEther now ranks above TRON.
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "'3.343.90"
$ws.Range("E18").Value = "  +1.90%  "
$ws.Range("D19").Value = "'17.90"
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D20").Value = "'0.117"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("E21").Value = "  -3.27%  "
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("D23").Value = "'18.43"
$ws.Range("E23").Value = "  +3.76%  "
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").Value = "'99.20"
$ws.Range("E25").Value = "  -7.17%  "
$ws.Range("E26").Value = "  -0.70%  "
$ws.Range("D27").Value = "'6.02"
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("E28").Value = "  +2.18%  "
$ws.Range("D29").Value = "'9.48"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "'8.55"
$ws.Range("E30").Value = "  -1.54%  "
$ws.Range("D31").Value = "'30.47"
$ws.Range("E31").Value = "  +0.75%  "
$ws.Range("D32").Value = "'6.70"
$ws.Range("E32").Value = "  +6.04%  "
$ws.Range("D33").Value = "'3.92"
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("D34").Value = "'592.57"
$ws.Range("E34").Value = "  +7.11%  "
$ws.Range("D35").Value = "'10.96"
$ws.Range("E35").Value = "  -0.80%  "
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").Value = "'3.700.79"
$ws.Range("E38").Value = "  +2.23%  "
$ws.Range("D39").Value = "'56.16"
$ws.Range("E39").Value = "  -2.04%  "
$ws.Range("D40").Value = "'3.50"
$ws.Range("E40").Value = "  -10.18%  "
$ws.Range("D41").Value = "'33.63"
$ws.Range("E41").Value = "  +5.25%  "
$ws.Range("D42").Value = "'0.0₃0702"
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("E44").Value = "  -8.68%  "
$ws.Range("E45").Value = "  -2.57%  "
$ws.Range("E46").Value = "  +4.40%  "
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").Value = "'2.56"
$ws.Range("E50").Value = "  -1.98%  "
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  +0.20%  "
